# New crime data collected - weekly CompStat refresh (30th Precinct)
# Moves the reporting window forward one week:
#   Volume 32 Number 12 -> 13
#   Report Covering the Week 3/17/2025 Through 3/23/2025
#     -> 3/24/2025 Through 3/30/2025
# and refreshes every Week-to-Date / 28-Day / Year-to-Date / historical
# percentage cell in the Crime Complaints table (rows 15-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header / report-window text -----------------------------------
$ws.Range("A8").Value = "Volume 32   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# ---- Crime Complaints table (rows 15-31) ----------------------------
# Numeric cells are written directly; cells that flip between the
# numeric-count style and the "0"/"***.*" placeholder-text style get a
# leading apostrophe (forces text) followed by a format copy/paste from
# an untouched donor cell that already carries the destination style
# (row 14: C14=text style, J14=integer style, K14=percent style).

$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("M15").Value = -75
$ws.Range("N15").Value = -88.235294117647

$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = -23.333333333333
$ws.Range("L16").Value = -8
$ws.Range("M16").Value = -45.238095238095
$ws.Range("N16").Value = -87.830687830687

$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 44
$ws.Range("K17").Value = 4.545454545454
$ws.Range("L17").Value = -11.538461538461
$ws.Range("M17").Value = 48.387096774193
$ws.Range("N17").Value = -72.781065088757

$ws.Range("D18").Value = 2
$ws.Range("J14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = -11.764705882352
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = -11.764705882352
$ws.Range("N18").Value = -92.268041237113

$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 175
$ws.Range("F19").Value = 24
$ws.Range("H19").Value = -4
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 89
$ws.Range("K19").Value = -17.977528089887
$ws.Range("L19").Value = 2.816901408450
$ws.Range("M19").Value = 121.212121212121
$ws.Range("N19").Value = -15.116279069767

$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 16
$ws.Range("K20").Value = -5.882352941176
$ws.Range("L20").Value = -38.461538461538
$ws.Range("M20").Value = 14.285714285714
$ws.Range("N20").Value = -77.142857142857

$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 166.666666666667
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 15.517241379310
$ws.Range("I21").Value = 175
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = -12.5
$ws.Range("L21").Value = -10.714285714285
$ws.Range("M21").Value = 19.863013698630
$ws.Range("N21").Value = -76.287262872628

$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = 100
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 16.666666666666
$ws.Range("L22").Value = -36.363636363636
$ws.Range("M22").Value = 133.333333333333

$ws.Range("D23").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$ws.Range("J14").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("H23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = -50

$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 72.727272727272
$ws.Range("F24").Value = 64
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 10.344827586206
$ws.Range("I24").Value = 181
$ws.Range("J24").Value = 190
$ws.Range("K24").Value = -4.736842105263
$ws.Range("L24").Value = 5.232558139534
$ws.Range("M24").Value = 135.064935064935

$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 38.461538461538
$ws.Range("I25").Value = 36
$ws.Range("J25").Value = 36
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = -14.285714285714

$ws.Range("C26").Value = 12
$ws.Range("E26").Value = 500
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 87
$ws.Range("J26").Value = 66
$ws.Range("K26").Value = 31.818181818181
$ws.Range("L26").Value = 24.285714285714
$ws.Range("M26").Value = -18.691588785046

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1

$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -57.142857142857

$ws.Range("N29").Value = -96.428571428571

$ws.Range("N30").Value = -96.296296296296

$ws.Range("D31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
